$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

# 1) Extend the "bioenergy" commodity into the pset_ci list on the ~TFM_INS
#    "start" row.
$ws.Range("F48").Value = "solar,wind,coal,gas,nuclear,hydro,bioenergy"

# 2) New ~TFM_INS row limiting LIFE for conventional/bioenergy processes,
#    appended right below the "start" row.
$ws.Range("C49").Value = "life"
$ws.Range("D49").Value = 40
$ws.Range("E49").Value = "'-life"
$ws.Range("F49").Value = "coal,gas,nuclear,bioenergy"

# 3) Split the combined wind RE mask/commodity into separate offshore and
#    onshore rows.
$ws.Range("C43").Value = "E[_]WOF*"
$ws.Range("D43").Value = "windoff"

$ws.Rows("44:44").Insert()
$ws.Range("C44").Value = "E[_]WON*"
$ws.Range("D44").Value = "windon"
$ws.Range("E44").Value = "IN"

# 4) ElcAgg_Wind input commodity mask: ELC_won* -> ELC_wo*
$ws.Range("D41").Value = "ELC_wo*"
